$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 86, pushing the existing rows 86-88 down to 87-89.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new record.
$ws.Cells.Item(86, 1).Value = 12
$ws.Cells.Item(86, 2).Value = "Mapocho Venta Directa de Santiago"
$ws.Cells.Item(86, 3).Value = "Metropolitana"
$ws.Cells.Item(86, 4).Value = 44826
$ws.Cells.Item(86, 4).Style = $ws.Cells.Item(87, 4).Style
$ws.Cells.Item(86, 4).NumberFormat = $ws.Cells.Item(87, 4).NumberFormat
$ws.Cells.Item(86, 5).Value = 13
$ws.Cells.Item(86, 6).Value = 100112002
$ws.Cells.Item(86, 7).Value = "Pimiento"
$ws.Cells.Item(86, 8).Value = "Zafiro rojo"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 80
$ws.Cells.Item(86, 11).Value = 15000
$ws.Cells.Item(86, 12).Value = 15000
$ws.Cells.Item(86, 13).Value = 15000
$ws.Cells.Item(86, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(86, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(86, 16).Value = 833
$ws.Cells.Item(86, 17).Value = 18
$ws.Cells.Item(86, 18).Value = "Hortaliza"
